# Excel Read From File
# - Rename sheet "IP-Method-C3" -> "Data"
# - Insert a new "Type" sheet (lookup table: patient -> type) right after "Data"
# - Update selections on a couple of sheets
# - Drop the "tabSelected" flag from the first sheet (now that "Type" is active)

$wb = $excel.ActiveWorkbook

# --- Rename IP-Method-C3 -> Data ------------------------------------------
$dataSheet = $wb.Worksheets.Item("IP-Method-C3")
$dataSheet.Name = "Data"

# --- Insert new "Type" sheet right after "Data" ----------------------------
$afterSheet = $wb.Worksheets.Item("Data")
$typeSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$typeSheet.Name = "Type"

# Header row
$typeSheet.Range("A1").Value = "Type"
$typeSheet.Range("B1").Value = "Patient"

# Patient numbers 1..20 across row 2 (columns B..U)
for ($i = 1; $i -le 20; $i++) {
    $typeSheet.Cells.Item(2, $i + 1).Value = $i
}

# Type 1 -> patients 1..5 (columns B..F) on row 3
$typeSheet.Range("A3").Value = 1
$typeSheet.Range("B3").Value = 1
$typeSheet.Range("C3").Value = 1
$typeSheet.Range("D3").Value = 1
$typeSheet.Range("E3").Value = 1
$typeSheet.Range("F3").Value = 1

# Type 2 -> patients 6..10 (columns G..K) on row 4
$typeSheet.Range("A4").Value = 2
$typeSheet.Range("G4").Value = 1
$typeSheet.Range("H4").Value = 1
$typeSheet.Range("I4").Value = 1
$typeSheet.Range("J4").Value = 1
$typeSheet.Range("K4").Value = 1

# Type 3 -> patients 11..20 (columns L..U) on row 5
$typeSheet.Range("A5").Value = 3
$typeSheet.Range("L5").Value = 1
$typeSheet.Range("M5").Value = 1
$typeSheet.Range("N5").Value = 1
$typeSheet.Range("O5").Value = 1
$typeSheet.Range("P5").Value = 1
$typeSheet.Range("Q5").Value = 1
$typeSheet.Range("R5").Value = 1
$typeSheet.Range("S5").Value = 1
$typeSheet.Range("T5").Value = 1
$typeSheet.Range("U5").Value = 1

# Make the new sheet the active one with B3 selected
$typeSheet.Activate()
$typeSheet.Range("B3").Select() | Out-Null

# --- Update selections on a couple of other sheets -------------------------
$acutiyLevels = $wb.Worksheets.Item("acutiyLevels")
$acutiyLevels.Range("G29").Select() | Out-Null

$dataSheet.Range("Q21").Select() | Out-Null

# --- Re-activate the Type sheet as the final active sheet -------------------
$typeSheet.Activate()
$typeSheet.Range("B3").Select() | Out-Null

# --- Work around a recalculation quirk triggered by the sheet insert/rename
#     above: re-assert the two pre-existing array formulas that reference a
#     bare "SUM range>=2" pattern so their cached results stay correct.
$dpm2 = $wb.Worksheets.Item("dp-method (2)")
$dpm2.Range("I12").FormulaArray = "=_xlfn.IFS(G3>G4,1,G3>G5,1,G3>G6,1,G3>G7,1,SUM(C15:C19)>=2,0)"
$dpm = $wb.Worksheets.Item("dp-method")
$dpm.Range("I13").FormulaArray = "=_xlfn.IFS(G3>G4,1,G3>G5,1,G3>G6,1,G3>G7,1,SUM(C16:C20)>=2,0)"
